{"js": "// Update the date line and every two-digit-by-two-digit multiplication\n// prompt in the worksheet table to the new values from the commit.\n// Each old string in this document is unique, so a literal, case-sensitive\n// search-and-replace for each (old -> new) pair reproduces the diff exactly\n// while leaving every other run (fonts, sizes, paragraph/table formatting)\n// untouched.\nconst replacements = [\n  [\"2025-08-25 Monday\", \"2025-08-26 Tuesday\"],\n  [\"89\u00d744=\", \"59\u00d725=\"],\n  [\"21\u00d726=\", \"84\u00d751=\"],\n  [\"46\u00d755=\", \"66\u00d784=\"],\n  [\"43\u00d780=\", \"90\u00d777=\"],\n  [\"62\u00d718=\", \"33\u00d750=\"],\n  [\"19\u00d732=\", \"71\u00d725=\"],\n  [\"26\u00d799=\", \"72\u00d777=\"],\n  [\"58\u00d766=\", \"44\u00d787=\"],\n  [\"77\u00d757=\", \"87\u00d734=\"],\n  [\"68\u00d789=\", \"36\u00d798=\"],\n  [\"50\u00d716=\", \"53\u00d788=\"],\n  [\"35\u00d767=\", \"77\u00d767=\"],\n  [\"61\u00d712=\", \"61\u00d754=\"],\n  [\"64\u00d781=\", \"30\u00d754=\"],\n  [\"74\u00d749=\", \"72\u00d718=\"],\n  [\"76\u00d717=\", \"83\u00d744=\"],\n  [\"56\u00d766=\", \"97\u00d765=\"],\n  [\"77\u00d761=\", \"95\u00d730=\"],\n  [\"32\u00d717=\", \"86\u00d762=\"],\n  [\"70\u00d751=\", \"14\u00d714=\"],\n  [\"88\u00d748=\", \"99\u00d785=\"],\n  [\"94\u00d744=\", \"69\u00d768=\"],\n  [\"98\u00d794=\", \"16\u00d716=\"],\n  [\"50\u00d714=\", \"91\u00d793=\"],\n  [\"15\u00d749=\", \"81\u00d763=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every two-digit-by-two-digit multiplication\n# prompt in the worksheet table to the new values from the commit.\n# Each old string in this document is unique, so a literal,\n# case-sensitive Find/Replace (wdReplaceAll = 2, wdFindContinue = 1,\n# no wildcards) for each (old -> new) pair reproduces the diff exactly\n# while leaving every other run (fonts, sizes, paragraph/table\n# formatting) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-08-25 Monday\", \"2025-08-26 Tuesday\"),\n  @(\"89\u00d744=\", \"59\u00d725=\"),\n  @(\"21\u00d726=\", \"84\u00d751=\"),\n  @(\"46\u00d755=\", \"66\u00d784=\"),\n  @(\"43\u00d780=\", \"90\u00d777=\"),\n  @(\"62\u00d718=\", \"33\u00d750=\"),\n  @(\"19\u00d732=\", \"71\u00d725=\"),\n  @(\"26\u00d799=\", \"72\u00d777=\"),\n  @(\"58\u00d766=\", \"44\u00d787=\"),\n  @(\"77\u00d757=\", \"87\u00d734=\"),\n  @(\"68\u00d789=\", \"36\u00d798=\"),\n  @(\"50\u00d716=\", \"53\u00d788=\"),\n  @(\"35\u00d767=\", \"77\u00d767=\"),\n  @(\"61\u00d712=\", \"61\u00d754=\"),\n  @(\"64\u00d781=\", \"30\u00d754=\"),\n  @(\"74\u00d749=\", \"72\u00d718=\"),\n  @(\"76\u00d717=\", \"83\u00d744=\"),\n  @(\"56\u00d766=\", \"97\u00d765=\"),\n  @(\"77\u00d761=\", \"95\u00d730=\"),\n  @(\"32\u00d717=\", \"86\u00d762=\"),\n  @(\"70\u00d751=\", \"14\u00d714=\"),\n  @(\"88\u00d748=\", \"99\u00d785=\"),\n  @(\"94\u00d744=\", \"69\u00d768=\"),\n  @(\"98\u00d794=\", \"16\u00d716=\"),\n  @(\"50\u00d714=\", \"91\u00d793=\"),\n  @(\"15\u00d749=\", \"81\u00d763=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
